# Update the "dSF" (column F) values for specific rows in Sheet1, as part of
# a data repull / mean-calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0
    3  = -4
    5  = -3
    7  = -4
    10 = -1
    11 = 6
    12 = -3
    16 = -3
    18 = -4
    19 = 5
    20 = 2
    21 = 3
    22 = -1
    25 = -4
    30 = 3
    36 = -3
    40 = -4
    42 = 0
    47 = -3
    53 = -1
    58 = -5
    61 = -6
    62 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
